# Added tc006 & tc007 of Toggle
#
# This edit:
#  1) Selects A1:B1 on the existing "CustomSDG" sheet (new selection there).
#  2) Duplicates the "CustomSDG" sheet (to inherit matching formatting/styles)
#     as a new sheet named "ToggleBtn", placed after it (last tab).
#  3) Trims the duplicated sheet down to just columns A:B, clears the
#     inherited formatting on row 2, and fills in the new toggle data.
#  4) Sets column B's width and the new sheet's selected cell.
#  5) Makes the new "ToggleBtn" sheet the active tab.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("CustomSDG")

# Update the selection on the CustomSDG sheet.
$src.Activate()
$src.Range("A1:B1").Select() | Out-Null

# Duplicate CustomSDG (preserves the header/row styles) to become the new
# ToggleBtn sheet, inserted as the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "ToggleBtn"

# Drop the extra columns C:D that CustomSDG had (ToggleBtn only needs A:B).
$newSheet.Range("C:D").EntireColumn.Delete()

# Row 2 on ToggleBtn should not carry over CustomSDG's row-2 styling.
$newSheet.Range("A2:B2").ClearFormats()

# Fill in the toggle button data (order chosen to match shared-string order).
$newSheet.Range("B2").Value = "Active Deals with All Stages !@#$%^&*() @#$%^&*Deals:Custom SDG"
$newSheet.Range("A2").Value = "TOGGLEBTN1"
$newSheet.Range("B1").Value = "Toggle_Button"

# Widen column B to fit the long toggle text.
$newSheet.Columns("B").ColumnWidth = 63.16666666666667

# Select E13 on the new sheet and make it the active sheet/tab.
$newSheet.Activate()
$newSheet.Range("E13").Select() | Out-Null
